$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the FilesTab Neo4j query text in cell B4 ---
# Removed the "File Type" and "Breed" RETURN lines from the cypher query
# (the following lines picked up an extra leading space as a result of the
# author's edit).
$ws.Range("B4").Value = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN['Golden Retriever']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
         coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
         coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

# The row shrinks (fewer wrapped lines) now that two RETURN clauses were
# dropped from the query text -- match the new auto-computed wrap height.
$ws.Rows.Item(4).RowHeight = 217.5

# --- Update the view/selection state ---
# Selection moves from the whole A column to cell B4, and the window
# scrolls so row 4 is the first visible row.
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 4
